$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '''58.895.73'
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = '  -0.60%  '
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = '''2.500.86'
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = '  +2.06%  '
$ws.Range("E3").Style = "Normal"

$ws.Range("E4").Value = '  +0.12%  '
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = '''540.31'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.06%  '
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = '''143.73'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -2.53%  '
$ws.Range("E6").Style = "Normal"

$ws.Range("E7").Value = '  -0.07%  '
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = '''0.572'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.56%  '
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = '''2.523.13'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.22%  '
$ws.Range("E9").Style = "Normal"

$ws.Range("E10").Value = '  +1.20%  '
$ws.Range("E10").Style = "Normal"

$ws.Range("E11").Value = '  +0.40%  '
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = '''5.59'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.61%  '
$ws.Range("E12").Style = "Normal"

$ws.Range("E13").Value = '  +1.33%  '
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = '''2.945.39'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +2.02%  '
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = '''23.55'
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -1.71%  '
$ws.Range("E15").Style = "Normal"

$ws.Range("D16").Value = '''58.834.67'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.67%  '
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = '''0.0000139'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +1.35%  '
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = '''2.524.99'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +0.43%  '
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = '''11.20'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  +0.79%  '
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = '  -1.56%  '
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = '''324.91'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.47%  '
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = '''0.997'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +3.08%  '
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = '''5.78'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.77%  '
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = '''61.98'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.50%  '
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = '''0.440'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -4.46%  '
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = '''0.162'
$ws.Range("D26").Style = "Normal"

$ws.Range("B27").Value = 'Binance-PegBSC-USD'
$ws.Range("B27").Style = "Normal"
$ws.Range("C27").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$ws.Range("C27").Style = "Normal"
$ws.Range("D27").Value = '''1.00'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +2.29%  '
$ws.Range("E27").Style = "Normal"

$ws.Range("B28").Value = 'WrappedeETH'
$ws.Range("B28").Style = "Normal"
$ws.Range("C28").Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range("C28").Style = "Normal"
$ws.Range("D28").Value = '''2.620.13'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +2.40%  '
$ws.Range("E28").Style = "Normal"

$ws.Range("E29").Value = '  +1.00%  '
$ws.Range("E29").Style = "Normal"

$ws.Range("B30").Value = 'PEPE'
$ws.Range("B30").Style = "Normal"
$ws.Range("C30").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("C30").Style = "Normal"
$ws.Range("D30").Value = '''0.0₃0773'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.66%  '
$ws.Range("E30").Style = "Normal"

$ws.Range("B31").Value = 'PancakeSwap'
$ws.Range("B31").Style = "Normal"
$ws.Range("C31").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("C31").Style = "Normal"
$ws.Range("D31").Value = '''1.81'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.16%  '
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = '''6.67'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  -1.86%  '
$ws.Range("E32").Style = "Normal"

$ws.Range("E33").Value = '  -4.52%  '
$ws.Range("E33").Style = "Normal"

$ws.Range("E34").Value = '  -0.12%  '
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = '''156.28'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +0.56%  '
$ws.Range("E35").Style = "Normal"

$ws.Range("E36").Value = '  +3.01%  '
$ws.Range("E36").Style = "Normal"

$ws.Range("D37").Value = '''18.67'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +1.68%  '
$ws.Range("E37").Style = "Normal"

$ws.Range("D38").Value = '''4.34'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -4.63%  '
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = '''1.58'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  -8.17%  '
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = '''5.69'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -3.11%  '
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = '''36.90'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +0.49%  '
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = '''294.98'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -5.63%  '
$ws.Range("E42").Style = "Normal"

$ws.Range("E43").Value = '  -0.30%  '
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = '''0.819'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -1.90%  '
$ws.Range("E44").Style = "Normal"

$ws.Range("E45").Value = '  -0.13%  '
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = '''0.600'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.52%  '
$ws.Range("E46").Style = "Normal"

$ws.Range("E47").Value = '  +0.46%  '
$ws.Range("E47").Style = "Normal"

$ws.Range("E48").Value = '  -0.88%  '
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = '''123.51'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.99%  '
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = '''18.55'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.40%  '
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = '  -2.18%  '
$ws.Range("E51").Style = "Normal"
